$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H53").Value = 353.70587
$ws.Range("I53").Value = 465.5
$ws.Range("K53").Value = 465.5
$ws.Range("M53").Value = 171.5

$ws.Range("H64").Value = 60389316
$ws.Range("I64").Value = 135870910
$ws.Range("J64").Value = 4038.6
$ws.Range("K64").Value = 135870910
$ws.Range("L64").Value = 4038.6
$ws.Range("M64").Value = -135870662
$ws.Range("N64").Value = -4534.6

$ws.Range("H67").Value = 60389316
$ws.Range("I67").Value = 135870910
$ws.Range("J67").Value = 4038.6
$ws.Range("K67").Value = 135870910
$ws.Range("L67").Value = 4038.6
$ws.Range("M67").Value = -135870052
$ws.Range("N67").Value = -5754.6

$ws.Range("H70").Value = 4425.8237
$ws.Range("I70").Value = 2964.1428
$ws.Range("J70").Value = 5449
$ws.Range("K70").Value = 8892.428400000001
$ws.Range("L70").Value = 16347
$ws.Range("M70").Value = -8622.428400000001
$ws.Range("N70").Value = -16887

$ws.Range("H73").Value = 4425.8237
$ws.Range("I73").Value = 2964.1428
$ws.Range("J73").Value = 5449
$ws.Range("K73").Value = 8892.428400000001
$ws.Range("L73").Value = 16347
$ws.Range("M73").Value = -7956.428400000001
$ws.Range("N73").Value = -18219

$ws.Range("H92").Value = 43479110
$ws.Range("I92").Value = 47619692
$ws.Range("J92").Value = 2999.5
$ws.Range("K92").Value = 47619692
$ws.Range("L92").Value = 2999.5
$ws.Range("M92").Value = -47618444
$ws.Range("N92").Value = -5495.5

$ws.Range("H137").Value = 5133.2
$ws.Range("I137").Value = 2459.2856
$ws.Range("J137").Value = 6173.0557
$ws.Range("K137").Value = 7377.8568
$ws.Range("L137").Value = 18519.1671
$ws.Range("M137").Value = -4827.8568
$ws.Range("N137").Value = -23619.1671

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5662.1953
$ws.Range("I61").Value = 5477.6577
$ws.Range("K61").Value = 5477.6577
$ws.Range("M61").Value = -5265.6577

$ws.Range("H63").Value = 4558.6
$ws.Range("I63").Value = 3400
$ws.Range("K63").Value = 3400
$ws.Range("M63").Value = -2714

$ws.Range("H66").Value = 4558.6
$ws.Range("I66").Value = 3400
$ws.Range("K66").Value = 17000
$ws.Range("M66").Value = -13568

$ws.Range("H74").Value = 15626636
$ws.Range("J74").Value = 1895.6
$ws.Range("L74").Value = 1895.6
$ws.Range("N74").Value = -3643.6

$ws.Range("H77").Value = 15626636
$ws.Range("J77").Value = 1895.6
$ws.Range("L77").Value = 9478
$ws.Range("N77").Value = -18214

$ws.Range("H122").Value = 4087.6897
$ws.Range("I122").Value = 2713.1428
$ws.Range("K122").Value = 8139.428400000001
$ws.Range("M122").Value = -5689.428400000001

$ws.Range("H132").Value = 17738.217
$ws.Range("J132").Value = 5744.222
$ws.Range("L132").Value = 17232.666
$ws.Range("N132").Value = -22292.666

$ws.Range("H136").Value = 5662.1953
$ws.Range("I136").Value = 5477.6577
$ws.Range("K136").Value = 16432.9731
$ws.Range("M136").Value = -13882.9731

$ws.Range("H139").Value = 75478.86
$ws.Range("J139").Value = 58897.332
$ws.Range("L139").Value = 58897.332
$ws.Range("N139").Value = -69177.33199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1157
$ws.Range("I134").Value = 1182.4445
$ws.Range("K134").Value = 3547.3335
$ws.Range("M134").Value = -1012.3335

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H19").Value = 61629.1
$ws.Range("I19").Value = 1786.5
$ws.Range("J19").Value = 300999.5
$ws.Range("K19").Value = 1786.5
$ws.Range("L19").Value = 300999.5
$ws.Range("M19").Value = -1616.5
$ws.Range("N19").Value = -301339.5

$ws.Range("H24").Value = 61629.1
$ws.Range("I24").Value = 1786.5
$ws.Range("J24").Value = 300999.5
$ws.Range("K24").Value = 1786.5
$ws.Range("L24").Value = 300999.5
$ws.Range("M24").Value = -1616.5
$ws.Range("N24").Value = -301339.5

$ws.Range("H31").Value = 18524270
$ws.Range("I31").Value = 71430616
$ws.Range("J31").Value = 7048.3
$ws.Range("K31").Value = 71430616
$ws.Range("L31").Value = 7048.3
$ws.Range("M31").Value = -71430321
$ws.Range("N31").Value = -7638.3

$ws.Range("H34").Value = 18524270
$ws.Range("I34").Value = 71430616
$ws.Range("J34").Value = 7048.3
$ws.Range("K34").Value = 71430616
$ws.Range("L34").Value = 7048.3
$ws.Range("M34").Value = -71430414
$ws.Range("N34").Value = -7452.3

$ws.Range("H58").Value = 197965.7
$ws.Range("I58").Value = 271403.25
$ws.Range("K58").Value = 271403.25
$ws.Range("M58").Value = -271200.25

$ws.Range("H62").Value = 62035.43
$ws.Range("I62").Value = 2750
$ws.Range("J62").Value = 71916.336
$ws.Range("K62").Value = 2750
$ws.Range("L62").Value = 71916.336
$ws.Range("M62").Value = -2126
$ws.Range("N62").Value = -73164.336

$ws.Range("H65").Value = 62035.43
$ws.Range("I65").Value = 2750
$ws.Range("J65").Value = 71916.336
$ws.Range("K65").Value = 13750
$ws.Range("L65").Value = 359581.68
$ws.Range("M65").Value = -10630
$ws.Range("N65").Value = -365821.68

$ws.Range("H132").Value = 55562296
$ws.Range("I132").Value = 78434430
$ws.Range("J132").Value = 15669.714
$ws.Range("K132").Value = 235303290
$ws.Range("L132").Value = 47009.142
$ws.Range("M132").Value = -235300760
$ws.Range("N132").Value = -52069.142

$ws.Range("H136").Value = 197965.7
$ws.Range("I136").Value = 271403.25
$ws.Range("K136").Value = 814209.75
$ws.Range("M136").Value = -811659.75

$ws.Range("H141").Value = 119152.43
$ws.Range("J141").Value = 119152.43
$ws.Range("L141").Value = 119152.43
$ws.Range("N141").Value = -129512.43

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H44").Value = 14286156
$ws.Range("I44").Value = 20000328
$ws.Range("J44").Value = 725
$ws.Range("K44").Value = 60000984
$ws.Range("L44").Value = 2175
$ws.Range("M44").Value = -60000586
$ws.Range("N44").Value = -2971

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H59").Value = 49000
$ws.Range("J59").Value = 49000
$ws.Range("L59").Value = 49000
$ws.Range("N59").Value = -50166

$ws.Range("H70").Value = 1628552.2
$ws.Range("I70").Value = 2530031.2
$ws.Range("K70").Value = 2530031.2
$ws.Range("M70").Value = -2529761.2

$ws.Range("H73").Value = 1628552.2
$ws.Range("I73").Value = 2530031.2
$ws.Range("K73").Value = 2530031.2
$ws.Range("M73").Value = -2529095.2

$ws.Range("H96").Value = 19980
$ws.Range("J96").Value = 19980
$ws.Range("L96").Value = 19980
$ws.Range("N96").Value = -25472

$ws.Range("H102").Value = 6409.8857
$ws.Range("I102").Value = 6355.154
$ws.Range("J102").Value = 6568
$ws.Range("K102").Value = 6355.154
$ws.Range("L102").Value = 6568
$ws.Range("M102").Value = -4733.154
$ws.Range("N102").Value = -9812

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6585848.5
$ws.Range("I62").Value = 7523827
$ws.Range("K62").Value = 7523827
$ws.Range("M62").Value = -7523203

$ws.Range("H65").Value = 6585848.5
$ws.Range("I65").Value = 7523827
$ws.Range("K65").Value = 37619135
$ws.Range("M65").Value = -37616015

$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()

$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()

$ws.Range("H122").Value = 4695.086
$ws.Range("I122").Value = 4211.3
$ws.Range("J122").Value = 7597.8
$ws.Range("K122").Value = 12633.9
$ws.Range("L122").Value = 22793.4
$ws.Range("M122").Value = -10183.9
$ws.Range("N122").Value = -27693.4

$ws.Range("H132").Value = 11629006
$ws.Range("I132").Value = 1006.30554
$ws.Range("J132").Value = 71430150
$ws.Range("K132").Value = 3018.91662
$ws.Range("L132").Value = 214290450
$ws.Range("M132").Value = -488.91662
$ws.Range("N132").Value = -214295510

$ws.Range("H136").Value = 5567.73
$ws.Range("I136").Value = 2228.0352
$ws.Range("J136").Value = 9994.768
$ws.Range("K136").Value = 6684.105599999999
$ws.Range("L136").Value = 29984.304
$ws.Range("M136").Value = -4134.105599999999
$ws.Range("N136").Value = -35084.304

